# UC6_TC1.xlsx - update evaluations (QuantitativeMetrics sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Assertion validity note (C7): re-worded, no longer references "baseline"
$ws.Range("C7").Value = "Test is not passing, but requirements says that  a message should be shown, but any message is shown."

# Code BLEU score (B12) and its breakdown note (C12) recomputed
$ws.Range("B12").Value = 0.2908983875687113
$ws.Range("C12").Value = "{'codebleu': 0.29089838756871134, 'ngram_match_score': 0.1658425703011601, 'weighted_ngram_match_score': 0.17588566580067874, 'syntax_match_score': 0.5266272189349113, 'dataflow_match_score': 0.29523809523809524}"

# Move/restore the active selection on the sheet to C7
[void]$ws.Activate()
[void]$ws.Range("C7").Select()
